# Build site at 2023-04-12 14:53:07 UTC
# Rewrites the "LOQ4069" course sheet: adds the Objectives/Teachers text
# that was missing, fixes the Programa resumido / Bibliografia content
# (which had been shifted one row off from their labels), and appends a
# new Requisitos content row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Rows 10-23 are being restructured (rows inserted/removed and
# heights changed), so drop the whole block and rebuild it from
# scratch to avoid carrying over stale row-height / cell metadata.
# ------------------------------------------------------------------
$ws.Range("A10:C23").EntireRow.Delete() | Out-Null
$ws.Range("A10:A24").EntireRow.Insert() | Out-Null

# Long multi-line strings, defined up front.
$objetivosTexto = @"
Complementar os conhecimentos na Área de Operações Unitárias da Indústria Química, com aplicações na operação, análise e projeto de equipamentos.
"@.Trim()

$programaResumidoTexto = @"
1. Tópicos especiais de operações unitárias envolvendo fluidos. 
2. Tópicos especiais de operações unitárias envolvendo transmissão de calor e massa.
"@.TrimEnd("`r","`n")

$metodoTexto = @"
O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios.
"@.Trim()

$bibliografiaTexto = @"
FOUST, Alan S. Princípios das Operações Unitárias. Rio de Janeiro : Guanabara Dois/LTC, 1982.
MCCABE, Warren; SMITH, Julian; HARRIOTT, Peter. Unit Operations of Chemical Engineering. Boston : McGraw-Hill, 2005.
GEANKOPLIS, Christie John. Transport Processes and Separation Process Principles. Upper Saddle River, NJ : Prentice Hall Professional Technical Reference, 2003.
COUPER, James R.; PENNEY, W. Roy; FAIR, James R.; WALAS, Stanley M. Chemical Process Equipment: Selection and Design. Amsterdam : Elsevier, c2005Boston.
PERRY, Robert H; GREEN, Don W; MALONEY, James O. Perry's Chemical Engineers' Handbook. 7th. ed. New York : McGraw-Hill, 1999.
Textos fornecidos pelo professor da disciplina
Artigos extraídos de revistas especializadas de Engenharia Química.
"@.Trim()

$requisitosTexto = "LOQ4057 -  Operações Unitárias III  (Requisito fraco)`n"

# Row 10 : Objetivos: / <new text> / <new text>
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = $objetivosTexto
$ws.Range("C10").Value = $objetivosTexto
$ws.Rows.Item(10).RowHeight = 60

# Row 11 : Objectives:  (B/C stay empty)
$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11:C11").Clear()
$ws.Rows.Item(11).RowHeight = 60

# Row 12 : Docentes responsáveis:  (B/C stay empty)
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12:C12").Clear()

# Row 13 : (no label, A stays empty) / 8151869 - Livia Chaguri e Carvalho
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C13").Value = "8151869 - Livia Chaguri e Carvalho"

# Row 14 : Programa resumido: / <new text> / <new text>
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = $programaResumidoTexto
$ws.Range("C14").Value = $programaResumidoTexto
$ws.Rows.Item(14).RowHeight = 60

# Row 15 : Short syllabus:  (B/C stay empty)
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15:C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# Row 16 : Programa: / <same text as Programa resumido>
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = $programaResumidoTexto
$ws.Range("C16").Value = $programaResumidoTexto
$ws.Rows.Item(16).RowHeight = 120

# Row 17 : Syllabus:  (B/C stay empty)
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17:C17").Clear()
$ws.Rows.Item(17).RowHeight = 120

# Row 18 : Avaliação:  (B/C stay empty)
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18:C18").Clear()

# Row 19 : Método: / <new text> / <new text>
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto
$ws.Rows.Item(19).RowHeight = 60

# Row 20 : Critério: / Provas e trabalhos.
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."
$ws.Rows.Item(20).RowHeight = 60

# Row 21 : Norma de recuperação: / Prova única com nota maior ou igual a 5,0 (cinco).
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Rows.Item(21).RowHeight = 60

# Row 22 : Bibliografia: / <new text> / <new text>
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = $bibliografiaTexto
$ws.Range("C22").Value = $bibliografiaTexto
$ws.Rows.Item(22).RowHeight = 120

# Row 23 : Requisitos:  (B/C stay empty)
$ws.Range("A23").Value = "Requisitos:"
$ws.Range("B23:C23").Clear()

# Row 24 (new) : (no label, A stays empty) / LOQ4057 requirement text
$ws.Range("A24").Clear()
$ws.Range("B24").Value = $requisitosTexto
$ws.Range("C24").Value = $requisitosTexto
$ws.Rows.Item(24).RowHeight = 30

# ------------------------------------------------------------------
# Column layout: column A alone keeps the label width/style; column B
# keeps its own (wider) width/style definition already in the sheet.
# Done last so it doesn't interfere with the cell writes above.
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 29.75
